$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove the old 4-column table entirely.
$ws.Cells.Clear()

# Table 1: Voertuig / eigenaar (rows 1-4)
$ws.Range("A1").Value = "[Voertuig]"
$ws.Range("B1").Value = "eigenaar"

$ws.Range("A2").Value = "Voertuig"
$ws.Range("B2").Value = "Persoon"

$ws.Range("A3").Value = "wagen1"
$ws.Range("B3").Value = "Piet"

$ws.Range("A4").Value = "boot2"
$ws.Range("B4").Value = "Marie"

# Row 5 intentionally left blank as a separator between the two tables.

# Table 2: Auto / aantalWielen (rows 6-8)
$ws.Range("A6").Value = "[Auto]"
$ws.Range("B6").Value = "aantalWielen"

$ws.Range("A7").Value = "Auto"

$ws.Range("A8").Value = "wagen1"
$ws.Range("B8").Value = 4
